# Weekly Fruit/Vegetable price update (Achicoria - Femacal de La Calera)
# Two new weekly rows are inserted into the existing data block:
#   - a new row at (current) row 117, pushing the former row 117 (and below) down by one
#   - a new row at (current) row 130 (after the first insert), pushing the former
#     rows 130/131 (now 131/132) down by one more
# The net effect matches the diff: dimension grows from A1:R131 to A1:R133 and
# two brand-new records are inserted in the middle of the table, with all
# subsequent rows shifted down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new rows -------------------------------------------------
$ws.Rows.Item(117).Insert()
$ws.Rows.Item(130).Insert()

# --- Fill in row 117 (new record, 2021-10-22) --------------------------------
$ws.Range("A117").Value2 = 3
$ws.Range("B117").Value2 = "Femacal de La Calera"
$ws.Range("C117").Value2 = "Coquimbo"
$ws.Range("D117").Value2 = 44491
$ws.Range("E117").Value2 = 5
$ws.Range("F117").Value2 = 100112010
$ws.Range("G117").Value2 = "Achicoria"
$ws.Range("H117").Value2 = "Sin especificar"
$ws.Range("I117").Value2 = "Primera"
$ws.Range("J117").Value2 = 120
$ws.Range("K117").Value2 = 5500
$ws.Range("L117").Value2 = 6000
$ws.Range("M117").Value2 = 5750
$ws.Range("N117").Value2 = "$/caja 16 unidades"
$ws.Range("O117").Value2 = "Provincia de Quillota"
$ws.Range("P117").Value2 = 359
$ws.Range("Q117").Value2 = 16
$ws.Range("R117").Value2 = "Hortaliza"

# --- Fill in row 130 (new record, 2021-09-22) --------------------------------
$ws.Range("A130").Value2 = 3
$ws.Range("B130").Value2 = "Femacal de La Calera"
$ws.Range("C130").Value2 = "Coquimbo"
$ws.Range("D130").Value2 = 44461
$ws.Range("E130").Value2 = 5
$ws.Range("F130").Value2 = 100112010
$ws.Range("G130").Value2 = "Achicoria"
$ws.Range("H130").Value2 = "Sin especificar"
$ws.Range("I130").Value2 = "Primera"
$ws.Range("J130").Value2 = 60
$ws.Range("K130").Value2 = 5000
$ws.Range("L130").Value2 = 5000
$ws.Range("M130").Value2 = 5000
$ws.Range("N130").Value2 = "$/caja 16 unidades"
$ws.Range("O130").Value2 = "Provincia de Quillota"
$ws.Range("P130").Value2 = 312
$ws.Range("Q130").Value2 = 16
$ws.Range("R130").Value2 = "Hortaliza"
